# Applies the "Merge in of data from master branch" edit.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Refresh the "_GoBack" bookmark to point at the very start of the
#    document. Word keeps only one "_GoBack" bookmark at a time, so
#    adding a new one removes the old one (which lived near the
#    k-means section) and renumbers every other bookmark in the
#    package, exactly as the target revision does (OLE_LINK1/2: 0,1 ->
#    1,2; OLE_LINK21/22: 2,3 -> 3,4; OLE_LINK3/4: 4,5 -> 5,6).
# ---------------------------------------------------------------------
$start = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $start) | Out-Null

# ---------------------------------------------------------------------
# 2. Update the course-code line on the title slide/page.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "COMP3005/Computer Vision", $true, $false, $false, $false, $false,
    $true, 1, $false, "COMP3204/COMP6223 - Computer Vision", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Page geometry: page height shrinks from 842pt (16840 twips) to
#    841pt (16820 twips); width is unchanged.
# ---------------------------------------------------------------------
$d.PageSetup.PageHeight = 841

# ---------------------------------------------------------------------
# 4. Footer: drop the stale "February 11, 2014" date stamp.
# ---------------------------------------------------------------------
$footers = $d.Sections(1).Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $f = $footers.Item($i)
    if ($f.Exists) {
        $f.Range.Find.Execute(
            "Written by Jonathon Hare, February 11, 2014.  Send any comments",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "Written by Jonathon Hare.  Send any comments", 2) | Out-Null
    }
}
